$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "LAI" column (AR) with header and per-row values
$ws.Range("AR1").Value = "LAI"
$ws.Range("AR2").Value = 3
$ws.Range("AR3").Value = 2
$ws.Range("AR4").Value = 2.5
$ws.Range("AR5").Value = 2
$ws.Range("AR6").Value = 2.5

# Replace the old "15.7+273" formulas in R4 and R5 with updated plain values
$ws.Range("R4").Value = 278.4
$ws.Range("R5").Value = 296.9

# Restore the final selection state left by the author
$ws.Range("H4").Select()
